# Update new cluster checklist to reflect the deprecation of the SSO server.
#
# The old "SSO" section (rows 15-20: SSO / server / app_id / app_secret /
# google client id / client secret) is replaced with a smaller
# "Authentication" section (rows 15-17: Authentication / Authentication
# provider + Google / OpenIDConnect / PAM / LDAP / credentials). This removes
# three rows, so everything below shifts up by three rows; the remaining
# section contents (S3 storage, Crunch credentials, Compute image, ssh
# access, Database password) stay the same, just three rows higher.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the three rows that are no longer needed (app_secret, google
# client id, client secret). This shifts all rows below up by 3, which
# automatically realigns the S3 storage / Crunch credentials / Compute
# image / ssh access / Database password sections to their new target
# rows while preserving their existing values and styles.
$ws.Rows.Item(18).Resize(3).Delete()

# Rewrite the (now 3-row) Authentication section.
$ws.Range("A15").Value = "Authentication"
$ws.Range("A16").Value = "Authentication provider"
$ws.Range("B16").Value = "Google / OpenIDConnect / PAM / LDAP"
$ws.Range("A17").Value = "credentials"
